$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.926.74'
$ws.Range('E2').Value = '  +2.95%  '
$ws.Range('D3').Value = '3.032.02'
$ws.Range('E3').Value = '  +1.84%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.38%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.028.15'
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.16%  '
$ws.Range('E11').Value = '  +4.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.464'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.13%  '
$ws.Range('E13').Value = '  +3.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.76%  '
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '3.534.07'
$ws.Range('E16').Value = '  +1.82%  '
$ws.Range('E17').Value = '  +2.72%  '
$ws.Range('D18').Value = '62.844.31'
$ws.Range('E18').Value = '  +2.73%  '
$ws.Range('D19').Value = '3.032.29'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '452.70'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.698'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.67%  '
$ws.Range('E23').Value = '  +3.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.91%  '
$ws.Range('E26').Value = '  +5.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.79%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.50'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.55%  '
$ws.Range('E30').Value = '  +10.58%  '
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.56'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('E34').Value = '  +1.84%  '
$ws.Range('D35').Value = '0.0₃0859'
$ws.Range('E35').Value = '  +5.43%  '
$ws.Range('E36').Value = '  +2.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.93'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.12%  '
$ws.Range('E38').Value = '  +12.75%  '
$ws.Range('E39').Value = '  +7.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.09'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.52'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.15'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.20%  '
$ws.Range('E43').Value = '  +15.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '44.22'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '390.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('E46').Value = '  +3.35%  '
$ws.Range('D47').Value = '2.722.05'
$ws.Range('E47').Value = '  +1.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.70'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('E50').Value = '  +7.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.85'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.08%  '
